# Commit: "added final permutation, convert key to ull"
#
# This adds a new "FinalPermutation" worksheet (built from a copy of the
# existing "Sheet1" -- which becomes "InitialPermutation" -- with the B/C
# columns swapped inside the formula text), reorders/renames the sheets,
# and nudges a couple of cell selections.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the Initial-Permutation sheet so the copy can become
#        the new Final-Permutation sheet, placed before the original. ---
$source = $wb.Worksheets.Item(1)
$source.Copy($source)

$finalPermutation = $wb.Worksheets.Item(1)
$initialPermutation = $wb.Worksheets.Item(2)
$sheet2 = $wb.Worksheets.Item(3)

$finalPermutation.Name = "FinalPermutation"
$initialPermutation.Name = "InitialPermutation"

# --- 2. Rewrite column A of FinalPermutation: swap the B/C references so
#        the table reads "Bit <dest> -> <src>" instead of "Bit <src> -> <dest>". ---
$finalPermutation.Range("A1").Formula = '="if (x & (i << "&64-C1&")) y |= (i << "&64-B1&"); //Bit "&C1&" -> "&B1'
$finalPermutation.Range("A2:A64").Formula = '="if (x & (i << "&64-C2&")) y |= (i << "&64-B2&"); //Bit "&C2&" -> "&B2'

# --- 3. View-state touch-ups (selections / active tab) ---
$finalPermutation.Activate()
$finalPermutation.Range("C1").Select()

$sheet2.Activate()
$sheet2.Range("T3").Select()

$initialPermutation.Activate()
$initialPermutation.Range("A38").Select()
